$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1, matching the style of the other header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data column F2:F22 -- time_taken metadata timestamps (no special style,
# matching the other plain data cells in columns B-E)
$timestamps = @(
    "2021-10-05 13:39:49.874054",
    "2021-10-05 13:39:49.874067",
    "2021-10-05 13:39:49.874071",
    "2021-10-05 13:39:49.874074",
    "2021-10-05 13:39:49.874078",
    "2021-10-05 13:39:49.874081",
    "2021-10-05 13:39:49.874084",
    "2021-10-05 13:39:49.874087",
    "2021-10-05 13:39:49.874090",
    "2021-10-05 13:39:49.874093",
    "2021-10-05 13:39:49.874096",
    "2021-10-05 13:39:49.874099",
    "2021-10-05 13:39:49.874102",
    "2021-10-05 13:39:49.874105",
    "2021-10-05 13:39:49.874108",
    "2021-10-05 13:39:49.874110",
    "2021-10-05 13:39:49.874114",
    "2021-10-05 13:39:49.874117",
    "2021-10-05 13:39:49.874120",
    "2021-10-05 13:39:49.874123",
    "2021-10-05 13:39:49.874126"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
